$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new kilométrage readings below the existing data (rows 36-37).
# Copy the existing date cell's format down first (same as dragging the fill
# handle) so the new cells reuse the same date style instead of creating a
# new one.
$ws.Range("A35").Copy($ws.Range("A36:A37"))

$ws.Range("A36").Value = 43787
$ws.Range("B36").Value = 1312

$ws.Range("A37").Value = 43795
$ws.Range("B37").Value = 1485

# Move the active selection to reflect the new first empty row, as Excel would
$ws.Range("B38").Select()

$wb.Save()
